# Rewrites the "34." / "35." changelog entries and appends new
# "36.", "37." and "38." entries, per the commit:
#   "Bridge and composite patterns implementations... replaced some methods..."
#
# Strategy: locate the three paragraphs that need to change (the
# "34. Introduced GameMessage..." paragraph, the "35. Removed all the
# rendering..." paragraph, and the paragraph that only carries the
# _GoBack bookmark), make room for three more paragraphs after the
# bookmark paragraph, then rewrite each of the six now-available
# paragraph slots in place with InsertXML so the exact run layout from
# the target OOXML is reproduced (including shared/non-shared runs).

$d = $word.ActiveDocument

$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.StartsWith("34. Introduced GameMessage")) {
        $startIndex = $i
        break
    }
}
if ($startIndex -eq -1) {
    throw "Could not locate target paragraph (34. Introduced GameMessage...)"
}

$i34 = $startIndex
$i35 = $startIndex + 1
$iBookmarkPara = $startIndex + 2

# Make room: insert two new empty paragraphs right before the
# bookmark-only paragraph (slots for the new "36." and "37." entries),
# and one new empty paragraph right after it (new trailing blank
# paragraph), keeping the bookmark paragraph itself as the slot for
# "38." so the _GoBack bookmark stays put.
$bp = $d.Paragraphs.Item($iBookmarkPara)

$before = $bp.Range.Duplicate()
$before.Collapse(1) | Out-Null
$before.InsertParagraphBefore()
$before.InsertParagraphBefore()

# Re-fetch: the bookmark paragraph index shifted down by 2.
$iBookmarkPara = $iBookmarkPara + 2
$bp = $d.Paragraphs.Item($iBookmarkPara)

$after = $bp.Range.Duplicate()
$after.Collapse(0) | Out-Null
$after.InsertParagraphAfter()

$i36 = $iBookmarkPara - 2
$i37 = $iBookmarkPara - 1
# ($iBookmarkPara + 1 is left as the new trailing empty paragraph.)

$xml34 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>3</w:t></w:r><w:r><w:t>4</w:t></w:r><w:r><w:t>. Removed all the rendering (printing) in the Renderer class.</w:t></w:r><w:r><w:t xml:space="preserve"> All renderable objects imp</w:t></w:r><w:r><w:t>lem</w:t></w:r><w:r><w:t>e</w:t></w:r><w:r><w:t xml:space="preserve">nt IRenderable interface and the Render method. </w:t></w:r><w:r><w:t xml:space="preserve">The game object should not now </w:t></w:r><w:r><w:t>how</w:t></w:r><w:r><w:t xml:space="preserve"> they will </w:t></w:r><w:r><w:t>be rendered for better coupling as well as th</w:t></w:r><w:r><w:t>e renderer does not care what will</w:t></w:r><w:r><w:t xml:space="preserve"> render.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>T</w:t></w:r><w:r><w:t xml:space="preserve">he game objects </w:t></w:r><w:r><w:t xml:space="preserve">can be rendered </w:t></w:r><w:r><w:t xml:space="preserve">with other </w:t></w:r><w:r><w:t>implementation</w:t></w:r><w:r><w:t xml:space="preserve"> easily </w:t></w:r><w:r><w:t>and objects can be added for rendering through the IRenderable interface.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml35 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>35</w:t></w:r><w:r><w:t>. Replaced Move method to Player class</w:t></w:r><w:r><w:t xml:space="preserve"> with the next cell checker altogether</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml36 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>36</w:t></w:r><w:r><w:t>. All renderable objects’ render methods receive particular IRenderer implementation as argument. It can be easily replaced with other implementation of the IRenderer. Bridge design pattern implementation.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml37 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>37</w:t></w:r><w:r><w:t>. Rendering the Maze renders each cell of it. Maze and the cells have tree-like structure.  Same with rendering the score list (each score item is rendered). Composite pattern implementation</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml38 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>38</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>. Maze solution checker method refactored to recursive solution for easier comprehension.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$d.Paragraphs.Item($i34).Range.InsertXML($xml34)
$d.Paragraphs.Item($i35).Range.InsertXML($xml35)
$d.Paragraphs.Item($i36).Range.InsertXML($xml36)
$d.Paragraphs.Item($i37).Range.InsertXML($xml37)
$d.Paragraphs.Item($iBookmarkPara).Range.InsertXML($xml38)

# Tidy the freshly-inserted trailing blank paragraph: InsertParagraphAfter
# leaves behind an empty-but-present run ("<w:r></w:r>"); replace it with
# a genuinely empty paragraph to match the surrounding blank paragraphs.
$iTrailingBlank = $iBookmarkPara + 1
$xmlEmpty = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item($iTrailingBlank).Range.InsertXML($xmlEmpty)
